# Updated cryptos list with refreshed price/volume figures.
# Column D ("Price") and column E ("Volume(1h)") values are plain text
# (they may contain thousands separators as literal dots, and percentage
# strings padded with spaces), so we force text formatting before writing
# each value to avoid Excel auto-converting number-like strings to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" '23.521.91'
Set-TextValue "E2" '  +1.43%  '
Set-TextValue "D3" '1.638.50'
Set-TextValue "E3" '  +2.31%  '
Set-TextValue "E4" '  +0.05%  '
Set-TextValue "D5" '307.66'
Set-TextValue "E5" '  +1.50%  '
Set-TextValue "E6" '  +0.04%  '
Set-TextValue "D7" '0.3765'
Set-TextValue "E7" '  -0.40%  '
Set-TextValue "D8" '52.53'
Set-TextValue "E8" '  +1.41%  '
Set-TextValue "D9" '0.3649'
Set-TextValue "E9" '  +1.01%  '
Set-TextValue "D10" '1.269'
Set-TextValue "E10" '  +0.50%  '
Set-TextValue "D11" '0.08177'
Set-TextValue "E11" '  +0.71%  '
Set-TextValue "D12" '1.002'
Set-TextValue "E12" '  +0.07%  '
Set-TextValue "D13" '23.00'
Set-TextValue "E13" '  +1.98%  '
Set-TextValue "D14" '6.637'
Set-TextValue "E15" '  +2.45%  '
Set-TextValue "D16" '7.401'
Set-TextValue "E16" '  +0.25%  '
Set-TextValue "D17" '1.639.05'
Set-TextValue "E17" '  +2.27%  '
Set-TextValue "D18" '94.66'
Set-TextValue "E18" '  +1.04%  '
Set-TextValue "D19" '0.06945'
Set-TextValue "E19" '  +1.22%  '
Set-TextValue "D20" '18.26'
Set-TextValue "E20" '  +1.40%  '
Set-TextValue "D21" '6.561'
Set-TextValue "E21" '  +0.50%  '
Set-TextValue "D22" '0.9999'
Set-TextValue "E22" '  -0.02%  '
Set-TextValue "D23" '23.516.25'
Set-TextValue "E23" '  +1.44%  '
Set-TextValue "E24" '  -0.96%  '
Set-TextValue "D25" '3.094'
Set-TextValue "E25" '  +3.74%  '
Set-TextValue "E26" '  +1.34%  '
Set-TextValue "D27" '21.29'
Set-TextValue "E27" '  +0.67%  '
Set-TextValue "D28" '151.58'
Set-TextValue "E28" '  +0.98%  '
Set-TextValue "D29" '5.347'
Set-TextValue "E29" '  +2.14%  '
Set-TextValue "D30" '135.45'
Set-TextValue "E30" '  +1.37%  '
Set-TextValue "D31" '2.377'
Set-TextValue "E31" '  -1.64%  '
Set-TextValue "D32" '1.818.56'
Set-TextValue "E32" '  +2.19%  '
Set-TextValue "D33" '6.802'
Set-TextValue "E33" '  -0.21%  '
Set-TextValue "D34" '0.9686'
Set-TextValue "E34" '  -0.83%  '
Set-TextValue "D35" '0.02826'
Set-TextValue "E35" '  +3.81%  '
Set-TextValue "D36" '10.31'
Set-TextValue "E36" '  -0.07%  '
Set-TextValue "D37" '0.07366'
Set-TextValue "E37" '  -2.09%  '
Set-TextValue "D38" '0.2544'
Set-TextValue "E38" '  +1.73%  '
Set-TextValue "D39" '6.181'
Set-TextValue "E39" '  +0.95%  '
Set-TextValue "E40" '  +0.95%  '
Set-TextValue "E41" '  +1.55%  '
Set-TextValue "D42" '0.7104'
Set-TextValue "E42" '  +0.13%  '
Set-TextValue "D43" '12.52'
Set-TextValue "E43" '  +0.61%  '
Set-TextValue "D44" '16.24'
Set-TextValue "E44" '  +5.18%  '
Set-TextValue "D45" '0.6546'
Set-TextValue "E45" '  +0.15%  '
Set-TextValue "D46" '2.342'
Set-TextValue "E46" '  +1.71%  '
Set-TextValue "D48" '4.043'
Set-TextValue "E48" '  +0.73%  '
Set-TextValue "D49" '0.07975'
Set-TextValue "E49" '  +0.29%  '
Set-TextValue "D50" '129.42'
Set-TextValue "E50" '  -2.13%  '
Set-TextValue "D51" '1.208'
Set-TextValue "E51" '  +0.45%  '
